$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated bio: VP of Pledge Education name was "TBD", now filled in.
$ws.Range("A4").Value = "Grace Till"
